$wb = $excel.ActiveWorkbook

# ALC row 8
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 127
$ws.Range("I8").Value = 156.25
$ws.Range("J8").Value = 10
$ws.Range("K8").Value = 468.75
$ws.Range("L8").Value = 30
$ws.Range("M8").Value = -329.75
$ws.Range("N8").Value = -308

# ARM row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2297.5625
$ws.Range("I2").Value = 2569.5833
$ws.Range("J2").Value = 1481.5
$ws.Range("K2").Value = 2569.5833
$ws.Range("L2").Value = 1481.5
$ws.Range("M2").Value = -2456.5833
$ws.Range("N2").Value = -1707.5

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 23205.96
$ws.Range("I32").Value = 24048.299
$ws.Range("K32").Value = 24048.299
$ws.Range("M32").Value = -23761.299

# ARM row 44
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 47087.5
$ws.Range("I44").Value = 10000
$ws.Range("J44").Value = 52385.715
$ws.Range("K44").Value = 10000
$ws.Range("L44").Value = 52385.715
$ws.Range("M44").Value = -9512
$ws.Range("N44").Value = -53361.715

# ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2072.35
$ws.Range("I45").Value = 2138.7222
$ws.Range("K45").Value = 2138.7222
$ws.Range("M45").Value = -1761.7222

# ARM row 55
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H55").Value = 69200
$ws.Range("J55").Value = 69200
$ws.Range("L55").Value = 69200
$ws.Range("N55").Value = -69830

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4787.6665
$ws.Range("I61").Value = 3488.361
$ws.Range("J61").Value = 8685.583000000001
$ws.Range("K61").Value = 3488.361
$ws.Range("L61").Value = 8685.583000000001
$ws.Range("M61").Value = -3276.361
$ws.Range("N61").Value = -9109.583000000001

# ARM row 97
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1310.9048
$ws.Range("I97").Value = 1246
$ws.Range("J97").Value = 1416.375
$ws.Range("K97").Value = 1246
$ws.Range("L97").Value = 1416.375
$ws.Range("M97").Value = -750
$ws.Range("N97").Value = -2408.375

# ARM row 110
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1652.4117
$ws.Range("I110").Value = 1652.7333
$ws.Range("J110").Value = 1650
$ws.Range("K110").Value = 1652.7333
$ws.Range("L110").Value = 1650
$ws.Range("M110").Value = 392.2666999999999
$ws.Range("N110").Value = -5740

# ARM row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 2297.5625
$ws.Range("I116").Value = 2569.5833
$ws.Range("J116").Value = 1481.5
$ws.Range("K116").Value = 2569.5833
$ws.Range("L116").Value = 1481.5
$ws.Range("M116").Value = -275.5832999999998
$ws.Range("N116").Value = -6069.5

# ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 4311864.5
$ws.Range("I122").Value = 1578.8462
$ws.Range("K122").Value = 4736.5386
$ws.Range("M122").Value = -2286.5386

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2191.3394
$ws.Range("I132").Value = 2075.575
$ws.Range("J132").Value = 2480.75
$ws.Range("K132").Value = 6226.724999999999
$ws.Range("L132").Value = 7442.25
$ws.Range("M132").Value = -3696.724999999999
$ws.Range("N132").Value = -12502.25

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 4787.6665
$ws.Range("I136").Value = 3488.361
$ws.Range("J136").Value = 8685.583000000001
$ws.Range("K136").Value = 10465.083
$ws.Range("L136").Value = 26056.749
$ws.Range("M136").Value = -7915.082999999999
$ws.Range("N136").Value = -31156.749

# BSM row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2297.5625
$ws.Range("I3").Value = 2569.5833
$ws.Range("J3").Value = 1481.5
$ws.Range("K3").Value = 2569.5833
$ws.Range("L3").Value = 1481.5
$ws.Range("M3").Value = -2455.5833
$ws.Range("N3").Value = -1709.5

# BSM row 20
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2516.7058
$ws.Range("I20").Value = 2985.3333
$ws.Range("J20").Value = 1392
$ws.Range("K20").Value = 2985.3333
$ws.Range("L20").Value = 1392
$ws.Range("M20").Value = -2738.3333
$ws.Range("N20").Value = -1886

# BSM row 80
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 293.1579
$ws.Range("J80").Value = 315.3125
$ws.Range("L80").Value = 315.3125
$ws.Range("N80").Value = -2311.3125

# BSM row 83
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H83").Value = 293.1579
$ws.Range("J83").Value = 315.3125
$ws.Range("L83").Value = 1576.5625
$ws.Range("N83").Value = -11560.5625

# BSM row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 8571.286
$ws.Range("I105").Value = 9199.799999999999
$ws.Range("K105").Value = 9199.799999999999
$ws.Range("M105").Value = -7452.799999999999

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2630.75
$ws.Range("I134").Value = 2732.8
$ws.Range("J134").Value = 1100
$ws.Range("K134").Value = 8198.400000000001
$ws.Range("L134").Value = 3300
$ws.Range("M134").Value = -5663.400000000001
$ws.Range("N134").Value = -8370

# CRP row 86
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 2659.8
$ws.Range("I86").Value = 2624.25
$ws.Range("K86").Value = 2624.25
$ws.Range("M86").Value = -1501.25

# CRP row 89
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 2659.8
$ws.Range("I89").Value = 2624.25
$ws.Range("K89").Value = 13121.25
$ws.Range("M89").Value = -7505.25

# CUL row 2
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 30.333334
$ws.Range("J2").Value = 33.48387
$ws.Range("L2").Value = 200.90322
$ws.Range("N2").Value = -426.90322

# CUL row 3
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 4239.1
$ws.Range("I3").Value = 2204.4443
$ws.Range("J3").Value = 7291.0835
$ws.Range("K3").Value = 6613.3329
$ws.Range("L3").Value = 21873.2505
$ws.Range("M3").Value = -6501.3329
$ws.Range("N3").Value = -22097.2505

# CUL row 7
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 220
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 220
$ws.Range("K7").Value = 0
$ws.Range("L7").ClearContents()
$ws.Range("M7").Value = 660
$ws.Range("N7").Value = -884

# CUL row 75
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 3951.0833
$ws.Range("I75").Value = 1978.25
$ws.Range("J75").Value = 4937.5
$ws.Range("K75").Value = 5934.75
$ws.Range("L75").Value = 14812.5
$ws.Range("M75").Value = -4936.75
$ws.Range("N75").Value = -16808.5

# CUL row 78
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H78").Value = 3951.0833
$ws.Range("I78").Value = 1978.25
$ws.Range("J78").Value = 4937.5
$ws.Range("K78").Value = 17804.25
$ws.Range("L78").Value = 44437.5
$ws.Range("M78").Value = -12812.25
$ws.Range("N78").Value = -54421.5

# CUL row 133
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H133").Value = 4179.8335
$ws.Range("I133").Value = 3110
$ws.Range("J133").Value = 5035.7
$ws.Range("K133").Value = 9330
$ws.Range("L133").Value = 15107.1
$ws.Range("M133").Value = -4270
$ws.Range("N133").Value = -25227.1

# GSM row 2
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 1129.9
$ws.Range("I2").Value = 99.85714
$ws.Range("J2").Value = 3533.3333
$ws.Range("K2").Value = 99.85714
$ws.Range("L2").Value = 3533.3333
$ws.Range("M2").Value = 13.14286
$ws.Range("N2").Value = -3759.3333

# GSM row 11
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 202916670
$ws.Range("I11").Value = 202916670
$ws.Range("K11").Value = 202916670
$ws.Range("M11").Value = -202916531

# GSM row 70
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5482.309
$ws.Range("I70").Value = 5356
$ws.Range("J70").Value = 5671.773
$ws.Range("K70").Value = 5356
$ws.Range("L70").Value = 5671.773
$ws.Range("M70").Value = -5086
$ws.Range("N70").Value = -6211.773

# GSM row 73
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 5482.309
$ws.Range("I73").Value = 5356
$ws.Range("J73").Value = 5671.773
$ws.Range("K73").Value = 5356
$ws.Range("L73").Value = 5671.773
$ws.Range("M73").Value = -4420
$ws.Range("N73").Value = -7543.773

# GSM row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 7578.615
$ws.Range("I80").Value = 15002.5
$ws.Range("J80").Value = 4279.1113
$ws.Range("K80").Value = 15002.5
$ws.Range("L80").Value = 4279.1113
$ws.Range("M80").Value = -14004.5
$ws.Range("N80").Value = -6275.1113

# GSM row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 7578.615
$ws.Range("I83").Value = 15002.5
$ws.Range("J83").Value = 4279.1113
$ws.Range("K83").Value = 75012.5
$ws.Range("L83").Value = 21395.5565
$ws.Range("M83").Value = -70020.5
$ws.Range("N83").Value = -31379.5565

# GSM row 97
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1980.5714
$ws.Range("I97").Value = 1980.5714
$ws.Range("K97").Value = 1980.5714
$ws.Range("M97").Value = -1484.5714

# GSM row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3320.25
$ws.Range("I102").Value = 2732.0715
$ws.Range("J102").Value = 3908.4285
$ws.Range("K102").Value = 2732.0715
$ws.Range("L102").Value = 3908.4285
$ws.Range("M102").Value = -1110.0715
$ws.Range("N102").Value = -7152.4285

# GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 17166.25
$ws.Range("I122").Value = 20555
$ws.Range("K122").Value = 61665
$ws.Range("M122").Value = -59215

# GSM row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2412.6785
$ws.Range("I126").Value = 1550
$ws.Range("J126").Value = 3275.3572
$ws.Range("K126").Value = 4650
$ws.Range("L126").Value = 9826.071599999999
$ws.Range("M126").Value = -2180
$ws.Range("N126").Value = -14766.0716

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 28686.764
$ws.Range("I132").Value = 43323.418
$ws.Range("J132").Value = 3595.3572
$ws.Range("K132").Value = 129970.254
$ws.Range("L132").Value = 10786.0716
$ws.Range("M132").Value = -127440.254
$ws.Range("N132").Value = -15846.0716

# LTW row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3279.5
$ws.Range("I7").Value = 3311.9375
$ws.Range("K7").Value = 3311.9375
$ws.Range("M7").Value = -3199.9375

# LTW row 23
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()

# LTW row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3124.64
$ws.Range("I40").Value = 2795.0557
$ws.Range("J40").Value = 3972.1428
$ws.Range("K40").Value = 2795.0557
$ws.Range("L40").Value = 3972.1428
$ws.Range("M40").Value = -2659.0557
$ws.Range("N40").Value = -4244.1428

# LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 6125.125
$ws.Range("I122").Value = 5868.5713
$ws.Range("K122").Value = 17605.7139
$ws.Range("M122").Value = -15155.7139

# LTW row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 3279.5
$ws.Range("I126").Value = 3311.9375
$ws.Range("K126").Value = 9935.8125
$ws.Range("M126").Value = -7465.8125

# WVR row 126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2006.25
$ws.Range("I126").Value = 1969.2307
$ws.Range("J126").Value = 2166.6667
$ws.Range("K126").Value = 5907.6921
$ws.Range("L126").Value = 6500.000100000001
$ws.Range("M126").Value = -3437.6921
$ws.Range("N126").Value = -11440.0001
